$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (F column = 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 42
$ws1.Range("F4").Value = 1546
$ws1.Range("F5").Value = 238
$ws1.Range("F7").Value = 932
$ws1.Range("F8").Value = 10065
$ws1.Range("F10").Value = 129
$ws1.Range("F11").Value = 248
$ws1.Range("F12").Value = 191
$ws1.Range("F13").Value = 383
$ws1.Range("F14").Value = 6978
$ws1.Range("F18").Value = 214

# Sheet "全部类型" updates (F column = 想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 42
$ws4.Range("F4").Value = 1546
$ws4.Range("F5").Value = 238
$ws4.Range("F8").Value = 932
$ws4.Range("F11").Value = 10065
$ws4.Range("F13").Value = 129
$ws4.Range("F14").Value = 248
$ws4.Range("F15").Value = 191
$ws4.Range("F16").Value = 383
$ws4.Range("F17").Value = 6978
$ws4.Range("F21").Value = 214
